$d = $word.ActiveDocument

$replacements = @(
    @("2023-11-28 Tuesday", "2023-11-29 Wednesday"),
    @("23×32=", "84×89="),
    @("41×28=", "83×95="),
    @("23×83=", "52×85="),
    @("41×76=", "36×80="),
    @("25×21=", "33×90="),
    @("68×14=", "26×88="),
    @("19×97=", "70×18="),
    @("77×95=", "48×23="),
    @("49×15=", "52×30="),
    @("75×74=", "66×22="),
    @("59×65=", "53×26="),
    @("42×75=", "29×28="),
    @("53×68=", "77×22="),
    @("72×94=", "44×94="),
    @("98×48=", "60×25="),
    @("99×48=", "83×93="),
    @("50×56=", "93×57="),
    @("49×67=", "86×29="),
    @("23×24=", "67×63="),
    @("64×25=", "94×17="),
    @("77×32=", "78×53="),
    @("71×78=", "99×24="),
    @("65×96=", "86×75="),
    @("59×27=", "93×64="),
    @("16×91=", "13×89=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
